$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the header formatting (bold, bordered, centered) used by the other header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New data values for columns I (I0) and J (IF)
$values = @(
    @(9, 9),
    @(9, 9),
    @(5, 5),
    @(6, 6),
    @(8, 8),
    @(9, 9),
    @(9, 9),
    @(8, 8),
    @(8, 8),
    @(9, 9),
    @(8, 8),
    @(9, 9),
    @(6, 6),
    @(8, 8),
    @(10, 10)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $pair = $values[$i]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
